$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "want to go" counts (column F)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 17
$wsExhibit.Range("F10").Value = 502

# Sheet "全部类型" (all types) - same underlying rows, mirror the updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1634
$wsAll.Range("F8").Value = 145
$wsAll.Range("F10").Value = 502
